# Feature: use cex_axis_y < 1 to make font smaller
# Insert two new settings rows (cex_axis_x, cex_axis_y) right after the
# existing "title_pt" row in the "# TEXT" section of base_settings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "# TEXT" section currently looks like:
#   row 11: # TEXT
#   row 12: decimal_sep | ,
#   row 13: title_pt    | 1.3 | times the default point size 7
#   row 14: (blank)
#   row 15: # PDF
#   ...
# We need to insert two new rows (14 and 15) for cex_axis_x / cex_axis_y,
# pushing everything from the old row 14 onward down by two rows.

$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = "cex_axis_x"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "times the default point size 7"

$ws.Range("A15").Value = "cex_axis_y"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "times the default point size 7"

# Reflect the author's final selection location.
$ws.Range("A15").Select()
